$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Date column (E2:E3) switches from real datetime values to plain text
#    timestamps ("dd-mm-yyyy hh:mm:ss" style strings). Write E3 first so the
#    shared-string table picks up the same ordering as the target workbook
#    (11:13:45 before 11:13:23).
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = "13-06-2024 11:13:45"
$ws.Range("E2").Value = "13-06-2024 11:13:23"

# ---------------------------------------------------------------------------
# 2. Header cell E1 gets a real date/time number format applied to it
#    (even though it holds text) -- mirrors the author re-using the old
#    date-format style slot for the header style instead.
# ---------------------------------------------------------------------------
$ws.Range("E1").NumberFormat = "dd/mm/yyyy\ hh:mm:ss"

# ---------------------------------------------------------------------------
# 3. New helper column F (print-alignment column) gets a light-gray thin
#    border on the right/top/bottom edges (no left edge) plus the default
#    Calibri 11 font, vertical-top + wrap alignment.
# ---------------------------------------------------------------------------
foreach ($r in 2, 3, 4) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Borders.Item(10).Color = 13882323
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Borders.Item(8).Color = 13882323
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).Color = 13882323
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Font.Name = "Calibri"
    $cell.VerticalAlignment = -4160
    $cell.WrapText = $true
}

# ---------------------------------------------------------------------------
# 4. New blank row 4 -- E4 continues the same bordered / Segoe UI styling as
#    E2:E3 (the "Date" column).
# ---------------------------------------------------------------------------
$e4 = $ws.Cells.Item(4, 5)
$e4.Borders.Item(7).Color = 13882323
$e4.Borders.Item(7).LineStyle = 1
$e4.Borders.Item(10).Color = 13882323
$e4.Borders.Item(10).LineStyle = 1
$e4.Borders.Item(8).Color = 13882323
$e4.Borders.Item(8).LineStyle = 1
$e4.Borders.Item(9).Color = 13882323
$e4.Borders.Item(9).LineStyle = 1
$e4.Font.Name = "Segoe UI"
$e4.Font.Size = 10
$e4.VerticalAlignment = -4160
$e4.WrapText = $true

# ---------------------------------------------------------------------------
# 5. Row heights / column width adjustments.
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 15.5
$ws.Rows.Item(4).RowHeight = 14.5
$ws.Columns.Item(5).ColumnWidth = 18.5703125

# ---------------------------------------------------------------------------
# 6. Selection + print setup (portrait page orientation).
# ---------------------------------------------------------------------------
$ws.Range("E7").Select()
$ws.PageSetup.Orientation = 1
